$wb = $excel.ActiveWorkbook

# This script applies a scheduled market-data refresh to the per-job profit
# columns (H:N) across all eight Disciple-of-the-Hand sheets. Values are raw
# Universalis price snapshots; nothing here is formula-driven.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 11917326
$ws.Range("I9").Value = 16667442
$ws.Range("J9").Value = 42034.5
$ws.Range("K9").Value = 16667442
$ws.Range("L9").Value = 42034.5
$ws.Range("M9").Value = -16667273
$ws.Range("N9").Value = -42372.5
$ws.Range("H28").Value = 266.41666
$ws.Range("I28").Value = 321.8889
$ws.Range("K28").Value = 321.8889
$ws.Range("M28").Value = 163.1111
$ws.Range("H58").Value = 608
$ws.Range("J58").Value = 900
$ws.Range("L58").Value = 2700
$ws.Range("N58").Value = -3000
$ws.Range("H74").Value = 11624.625
$ws.Range("I74").Value = 11624.625
$ws.Range("K74").Value = 11624.625
$ws.Range("M74").Value = -10688.625
$ws.Range("H77").Value = 11624.625
$ws.Range("I77").Value = 11624.625
$ws.Range("K77").Value = 58123.125
$ws.Range("M77").Value = -53443.125
$ws.Range("H80").Value = 46297304
$ws.Range("I80").Value = 83334056
$ws.Range("K80").Value = 250002168
$ws.Range("M80").Value = -250001170
$ws.Range("H83").Value = 46297304
$ws.Range("I83").Value = 83334056
$ws.Range("K83").Value = 750006504
$ws.Range("M83").Value = -750001512
$ws.Range("H125").Value = 252841000
$ws.Range("J125").Value = 320
$ws.Range("L125").Value = 2880
$ws.Range("N125").Value = -7800
$ws.Range("H132").Value = 1138.3611
$ws.Range("I132").Value = 969.4545000000001
$ws.Range("K132").Value = 2908.3635
$ws.Range("M132").Value = -378.3635000000004
$ws.Range("H141").Value = 1717.8572
$ws.Range("I141").Value = 1688.1578
$ws.Range("K141").Value = 5064.4734
$ws.Range("M141").Value = 115.5266000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 555.4761999999999
$ws.Range("I5").Value = 656.6923
$ws.Range("K5").Value = 656.6923
$ws.Range("M5").Value = -544.6923
$ws.Range("H32").Value = 259418
$ws.Range("I32").Value = 313755.7
$ws.Range("K32").Value = 313755.7
$ws.Range("M32").Value = -313468.7
$ws.Range("H45").Value = 44191.875
$ws.Range("I45").Value = 64939.312
$ws.Range("J45").Value = 2697
$ws.Range("K45").Value = 64939.312
$ws.Range("L45").Value = 2697
$ws.Range("M45").Value = -64562.312
$ws.Range("N45").Value = -3451
$ws.Range("H61").Value = 728575.4
$ws.Range("I61").Value = 2083.0322
$ws.Range("K61").Value = 2083.0322
$ws.Range("M61").Value = -1871.0322
$ws.Range("H97").Value = 10765.6
$ws.Range("I97").Value = 13032
$ws.Range("K97").Value = 13032
$ws.Range("M97").Value = -12536
$ws.Range("H122").Value = 2959.2
$ws.Range("I122").Value = 2959.2
$ws.Range("K122").Value = 8877.599999999999
$ws.Range("M122").Value = -6427.599999999999
$ws.Range("H130").Value = 49750
$ws.Range("J130").Value = 49750
$ws.Range("L130").Value = 49750
$ws.Range("N130").Value = -59790
$ws.Range("H132").Value = 1844.6753
$ws.Range("I132").Value = 1674.5211
$ws.Range("K132").Value = 5023.5633
$ws.Range("M132").Value = -2493.5633
$ws.Range("H135").Value = 92160.875
$ws.Range("J135").Value = 93183.86
$ws.Range("L135").Value = 93183.86
$ws.Range("N135").Value = -103323.86
$ws.Range("H136").Value = 728575.4
$ws.Range("I136").Value = 2083.0322
$ws.Range("K136").Value = 6249.096600000001
$ws.Range("M136").Value = -3699.096600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 555.4761999999999
$ws.Range("I4").Value = 656.6923
$ws.Range("K4").Value = 656.6923
$ws.Range("M4").Value = -541.6923
$ws.Range("H86").Value = 2806.923
$ws.Range("I86").Value = 1732.4445
$ws.Range("J86").Value = 5224.5
$ws.Range("K86").Value = 1732.4445
$ws.Range("L86").Value = 5224.5
$ws.Range("M86").Value = -609.4445000000001
$ws.Range("N86").Value = -7470.5
$ws.Range("H89").Value = 2806.923
$ws.Range("I89").Value = 1732.4445
$ws.Range("J89").Value = 5224.5
$ws.Range("K89").Value = 8662.2225
$ws.Range("L89").Value = 26122.5
$ws.Range("M89").Value = -3046.2225
$ws.Range("N89").Value = -37354.5
$ws.Range("H105").Value = 10435.857
$ws.Range("J105").Value = 5167.5
$ws.Range("L105").Value = 5167.5
$ws.Range("N105").Value = -8661.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3201.8975
$ws.Range("I31").Value = 2973.5908
$ws.Range("J31").Value = 3497.353
$ws.Range("K31").Value = 2973.5908
$ws.Range("L31").Value = 3497.353
$ws.Range("M31").Value = -2678.5908
$ws.Range("N31").Value = -4087.353
$ws.Range("H34").Value = 3201.8975
$ws.Range("I34").Value = 2973.5908
$ws.Range("J34").Value = 3497.353
$ws.Range("K34").Value = 2973.5908
$ws.Range("L34").Value = 3497.353
$ws.Range("M34").Value = -2771.5908
$ws.Range("N34").Value = -3901.353
$ws.Range("H122").Value = 2798.3572
$ws.Range("I122").Value = 2798.3572
$ws.Range("K122").Value = 8395.071599999999
$ws.Range("M122").Value = -5945.071599999999
$ws.Range("H132").Value = 20514.793
$ws.Range("I132").Value = 22279.062
$ws.Range("J132").Value = 3577.8
$ws.Range("K132").Value = 66837.186
$ws.Range("L132").Value = 10733.4
$ws.Range("M132").Value = -64307.186
$ws.Range("N132").Value = -15793.4
$ws.Range("H134").Value = 1862.2963
$ws.Range("I134").Value = 1794.8
$ws.Range("J134").Value = 1946.6666
$ws.Range("K134").Value = 5384.4
$ws.Range("L134").Value = 5839.9998
$ws.Range("M134").Value = -2849.4
$ws.Range("N134").Value = -10909.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2090.1428
$ws.Range("I121").Value = 547.6667
$ws.Range("K121").Value = 1643.0001
$ws.Range("M121").Value = -333.0001
$ws.Range("H137").Value = 3991.4167
$ws.Range("J137").Value = 3797
$ws.Range("L137").Value = 11391
$ws.Range("N137").Value = -21591

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H107").Value = 63217.625
$ws.Range("I107").Value = 200303.6
$ws.Range("J107").Value = 905.8182
$ws.Range("K107").Value = 200303.6
$ws.Range("L107").Value = 905.8182
$ws.Range("M107").Value = -198383.6
$ws.Range("N107").Value = -4745.8182
$ws.Range("H122").Value = 2209.28
$ws.Range("I122").Value = 2042.6364
$ws.Range("K122").Value = 6127.9092
$ws.Range("M122").Value = -3677.9092
$ws.Range("H126").Value = 2165
$ws.Range("J126").Value = 2000
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 662360.6
$ws.Range("I132").Value = 5688.5386
$ws.Range("K132").Value = 17065.6158
$ws.Range("M132").Value = -14535.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1016.7273
$ws.Range("J16").Value = 984.7143
$ws.Range("L16").Value = 984.7143
$ws.Range("N16").Value = -1324.7143
$ws.Range("H40").Value = 3536.2
$ws.Range("I40").Value = 1620.5
$ws.Range("K40").Value = 1620.5
$ws.Range("M40").Value = -1484.5
$ws.Range("H55").Value = 1155.1111
$ws.Range("I55").Value = 876.4
$ws.Range("J55").Value = 1503.5
$ws.Range("K55").Value = 876.4
$ws.Range("L55").Value = 1503.5
$ws.Range("M55").Value = -703.4
$ws.Range("N55").Value = -1849.5
$ws.Range("H61").Value = 2611.926
$ws.Range("I61").Value = 2417.625
$ws.Range("J61").Value = 4166.3335
$ws.Range("K61").Value = 2417.625
$ws.Range("L61").Value = 4166.3335
$ws.Range("M61").Value = -2215.625
$ws.Range("N61").Value = -4570.3335
$ws.Range("H82").Value = 832.94116
$ws.Range("I82").Value = 651.7143
$ws.Range("J82").Value = 959.8
$ws.Range("K82").Value = 651.7143
$ws.Range("L82").Value = 959.8
$ws.Range("M82").Value = -290.7143
$ws.Range("N82").Value = -1681.8
$ws.Range("H85").Value = 832.94116
$ws.Range("I85").Value = 651.7143
$ws.Range("J85").Value = 959.8
$ws.Range("K85").Value = 651.7143
$ws.Range("L85").Value = 959.8
$ws.Range("M85").Value = 596.2857
$ws.Range("N85").Value = -3455.8
$ws.Range("H113").Value = 2611.926
$ws.Range("I113").Value = 2417.625
$ws.Range("J113").Value = 4166.3335
$ws.Range("K113").Value = 2417.625
$ws.Range("L113").Value = 4166.3335
$ws.Range("M113").Value = -247.625
$ws.Range("N113").Value = -8506.333500000001
$ws.Range("H122").Value = 3534.7693
$ws.Range("I122").Value = 3095.1904
$ws.Range("J122").Value = 4047.611
$ws.Range("K122").Value = 9285.5712
$ws.Range("L122").Value = 12142.833
$ws.Range("M122").Value = -6835.5712
$ws.Range("N122").Value = -17042.833
$ws.Range("H132").Value = 2652.465
$ws.Range("I132").Value = 2389.258
$ws.Range("K132").Value = 7167.773999999999
$ws.Range("M132").Value = -4637.773999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 47328
$ws.Range("J46").Value = 47328
$ws.Range("L46").Value = 47328
$ws.Range("N46").Value = -47790
$ws.Range("H107").Value = 866414.2
$ws.Range("I107").Value = 630.4483
$ws.Range("K107").Value = 1891.3449
$ws.Range("M107").Value = 28.65509999999995
$ws.Range("H122").Value = 1733
$ws.Range("I122").Value = 1472.8214
$ws.Range("J122").Value = 3190
$ws.Range("K122").Value = 4418.4642
$ws.Range("L122").Value = 9570
$ws.Range("M122").Value = -1968.4642
$ws.Range("N122").Value = -14470
$ws.Range("H134").Value = 47328
$ws.Range("J134").Value = 47328
$ws.Range("L134").Value = 141984
$ws.Range("N134").Value = -147054
